$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

function Set-PlainCell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 2
Set-PlainCell "D2" "27.674.78"
Set-PlainCell "E2" "  +1.10%  "

# Row 3
Set-PlainCell "D3" "1.874.04"
Set-PlainCell "E3" "  +0.83%  "

# Row 4
Set-PlainCell "E4" "  +0.24%  "

# Row 5
Set-TextCell "D5" "331.61"
Set-PlainCell "E5" "  +2.52%  "

# Row 6
Set-TextCell "D6" "1.004"
Set-PlainCell "E6" "  +0.20%  "

# Row 7
Set-TextCell "D7" "0.4719"
Set-PlainCell "E7" "  +4.20%  "

# Row 8
Set-TextCell "D8" "0.3943"
Set-PlainCell "E8" "  +2.02%  "

# Row 9
Set-TextCell "D9" "47.92"
Set-PlainCell "E9" "  -1.43%  "

# Row 10
Set-TextCell "D10" "0.08032"
Set-PlainCell "E10" "  +1.55%  "

# Row 11
Set-PlainCell "E11" "  +1.04%  "

# Row 12
Set-TextCell "D12" "22.03"
Set-PlainCell "E12" "  +3.08%  "

# Row 13
Set-PlainCell "D13" "1.852.60"
Set-PlainCell "E13" "  -1.05%  "

# Row 14
Set-TextCell "D14" "5.965"
Set-PlainCell "E14" "  +0.83%  "

# Row 15
Set-TextCell "D15" "7.118"
Set-PlainCell "E15" "  +0.01%  "

# Row 16
Set-TextCell "D16" "1.004"
Set-PlainCell "E16" "  +0.27%  "

# Row 17
Set-TextCell "D17" "0.00001049"
Set-PlainCell "E17" "  +1.60%  "

# Row 18
Set-TextCell "D18" "87.02"
Set-PlainCell "E18" "  +1.41%  "

# Row 19
Set-TextCell "D19" "0.06684"
Set-PlainCell "E19" "  +2.50%  "

# Row 20
Set-TextCell "D20" "17.17"
Set-PlainCell "E20" "  +0.61%  "

# Row 21
Set-PlainCell "E21" "  +0.20%  "

# Row 22
Set-PlainCell "D22" "27.679.25"
Set-PlainCell "E22" "  +1.11%  "

# Row 23
Set-TextCell "D23" "5.516"
Set-PlainCell "E23" "  -0.23%  "

# Row 24
Set-TextCell "D24" "10.97"
Set-PlainCell "E24" "  +1.31%  "

# Row 25
Set-TextCell "D25" "2.308"
Set-PlainCell "E25" "  +1.22%  "

# Row 26
Set-PlainCell "D26" "2.087.19"
Set-PlainCell "E26" "  -0.11%  "

# Row 27
Set-TextCell "D27" "158.46"
Set-PlainCell "E27" "  +3.16%  "

# Row 28
Set-TextCell "D28" "20.15"
Set-PlainCell "E28" "  +2.19%  "

# Row 29
Set-TextCell "D29" "2.102"
Set-PlainCell "E29" "  +1.53%  "

# Row 30
Set-TextCell "D30" "5.576"
Set-PlainCell "E30" "  +2.50%  "

# Row 31
Set-TextCell "D31" "122.19"
Set-PlainCell "E31" "  +1.08%  "

# Row 32
Set-TextCell "D32" "0.9745"
Set-PlainCell "E32" "  +4.23%  "

# Row 33
Set-TextCell "D33" "0.09532"
Set-PlainCell "E33" "  +2.54%  "

# Row 34
Set-TextCell "D34" "1.446"
Set-PlainCell "E34" "  -2.73%  "

# Row 35
Set-TextCell "D35" "3.593"
Set-PlainCell "E35" "  -0.20%  "

# Row 36
Set-TextCell "D36" "5.335"
Set-PlainCell "E36" "  +1.65%  "

# Row 37
Set-TextCell "D37" "0.06101"
Set-PlainCell "E37" "  +1.85%  "

# Row 38
Set-TextCell "D38" "0.02251"
Set-PlainCell "E38" "  +0.77%  "

# Row 39
Set-TextCell "D39" "1.230"
Set-PlainCell "E39" "  +0.56%  "

# Row 40
Set-TextCell "D40" "8.231"
Set-PlainCell "E40" "  +0.18%  "

# Row 41
Set-TextCell "D41" "0.6023"
Set-PlainCell "E41" "  +2.08%  "

# Row 42
Set-TextCell "D42" "0.1905"
Set-PlainCell "E42" "  +0.93%  "

# Row 43
Set-PlainCell "E43" "  +1.49%  "

# Row 44
Set-TextCell "D44" "1.269"
Set-PlainCell "E44" "  -0.77%  "

# Row 45
Set-TextCell "D45" "0.5691"
Set-PlainCell "E45" "  +1.34%  "

# Row 46
Set-TextCell "D46" "12.17"
Set-PlainCell "E46" "  +1.46%  "

# Row 47
Set-TextCell "D47" "1.944"
Set-PlainCell "E47" "  +1.21%  "

# Row 48
Set-TextCell "D48" "3.379"
Set-PlainCell "E48" "  +0.26%  "

# Row 49
Set-PlainCell "B49" "Cronos"
Set-PlainCell "C49" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D49" "0.06888"
Set-PlainCell "E49" "  +1.76%  "

# Row 50
Set-PlainCell "B50" "Quant"
Set-PlainCell "C50" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D50" "115.14"
Set-PlainCell "E50" "  +6.49%  "

# Row 51
Set-TextCell "D51" "0.00000000300"
Set-PlainCell "E51" "  +9.29%  "
